$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values updated per the latest cryptos snapshot.
# Numeric-looking price strings (single-dot decimals) need to be forced
# back to text so Excel does not silently convert them to numbers -
# matching the original inlineStr/text storage of column D.

$ws.Range('D2').Value = '27.535.14'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '1.595.28'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.500'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.91%  '
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.25'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.04%  '
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0869'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = '1.823.12'
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').Value = '1.619.44'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('E14').Value = '  -3.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.539'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.32%  '
$ws.Range('D17').Value = '27.523.18'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '216.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.32%  '
$ws.Range('D20').Value = '0.0₃0690'
$ws.Range('E20').Value = '  -3.71%  '
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('E22').Value = '  -2.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.39%  '
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.64'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('E28').Value = '  -2.78%  '
$ws.Range('E29').Value = '  -4.60%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0467'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.28%  '
$ws.Range('E32').Value = '  -2.39%  '
$ws.Range('D33').Value = '1.358.62'
$ws.Range('E33').Value = '  -1.95%  '
$ws.Range('E34').Value = '  -3.66%  '
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.959'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.537'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.813'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.13%  '
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.965'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.85%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.32'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.84'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('E45').Value = '  -3.22%  '
$ws.Range('D46').Value = '1.733.33'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('E47').Value = '  -3.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.27%  '
$ws.Range('D49').Value = '0.0₇0998'
$ws.Range('E49').Value = '  +3.16%  '
$ws.Range('E50').Value = '  -3.47%  '
$ws.Range('E51').Value = '  -0.91%  '
